# Added Test Scripts: dealAnalysis_Pricing_263023_TC_01, dealAnalysis_Pricing_263025_TC_03
# for Deal Pricing Division (plus quote_ProposalSetup_262997_TC_01), appended to the
# "summary" sheet's test-results table (columns: S.No, TestCaseID, Testcase Name,
# Execution Date, Test Reult, Comments).

$wb = $excel.ActiveWorkbook

$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "summary") {
        $ws = $sheet
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

# New rows appended after the existing last row (15.0 at row 16).
$newRows = @(
    @{ SNo = 16; TestCaseID = "quote_ProposalSetup_262997_TC_01"; ExecDate = "25/05/2022"; Result = "Pass" },
    @{ SNo = 17; TestCaseID = "dealAnalysis_Pricing_263023_TC_01"; ExecDate = "27/05/2022"; Result = "Pass" },
    @{ SNo = 18; TestCaseID = "dealAnalysis_Pricing_263025_TC_03"; ExecDate = "27/05/2022"; Result = "Pass" }
)

$startRow = 17
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row.SNo        # A: S.No
    $ws.Cells.Item($r, 2).Value = $row.TestCaseID # B: TestCaseID
    $ws.Cells.Item($r, 4).Value = $row.ExecDate   # D: Execution Date
    $ws.Cells.Item($r, 5).Value = $row.Result     # E: Test Reult
}
